$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6244955548078103
$ws.Cells.Item(2, 3).Value = 0.4810815906795829
$ws.Cells.Item(2, 4).Value = 0.07932307540037664
$ws.Cells.Item(2, 5).Value = 0.211491955676601
$ws.Cells.Item(2, 6).Value = 0.4155727922916412
$ws.Cells.Item(2, 7).Value = 0.2461107820272446
$ws.Cells.Item(2, 8).Value = 2.234782218933105
$ws.Cells.Item(2, 9).Value = 1.181956171989441

$ws.Cells.Item(3, 2).Value = 0.6280056471781839
$ws.Cells.Item(3, 3).Value = 0.4898282375023925
$ws.Cells.Item(3, 4).Value = 0.09570534380868845
$ws.Cells.Item(3, 5).Value = 0.2254411426725104
$ws.Cells.Item(3, 6).Value = 0.4116881191730499
$ws.Cells.Item(3, 7).Value = 0.2419624626636505
$ws.Cells.Item(3, 8).Value = 2.195017099380493
$ws.Cells.Item(3, 9).Value = 1.161046504974365

$ws.Cells.Item(4, 2).Value = 0.6308014574709552
$ws.Cells.Item(4, 3).Value = 0.4994507942384036
$ws.Cells.Item(4, 4).Value = 0.1099278179446547
$ws.Cells.Item(4, 5).Value = 0.2378903412238682
$ws.Cells.Item(4, 6).Value = 0.4085939824581146
$ws.Cells.Item(4, 7).Value = 0.237398698925972
$ws.Cells.Item(4, 8).Value = 2.160494565963745
$ws.Cells.Item(4, 9).Value = 1.142385601997375

$ws.Cells.Item(5, 2).Value = 0.6338880409953453
$ws.Cells.Item(5, 3).Value = 0.5107600691727181
$ws.Cells.Item(5, 4).Value = 0.1247942637020881
$ws.Cells.Item(5, 5).Value = 0.2511136390873432
$ws.Cells.Item(5, 6).Value = 0.4051780700683594
$ws.Cells.Item(5, 7).Value = 0.2320349961519241
$ws.Cells.Item(5, 8).Value = 2.124408721923828
$ws.Cells.Item(5, 9).Value = 1.122564196586609

$ws.Cells.Item(6, 2).Value = 0.6361135386749246
$ws.Cells.Item(6, 3).Value = 0.5190702446017985
$ws.Cells.Item(6, 4).Value = 0.1376702844729901
$ws.Cells.Item(6, 5).Value = 0.2623179205384651
$ws.Cells.Item(6, 6).Value = 0.4027151465415955
$ws.Cells.Item(6, 7).Value = 0.2280936539173126
$ws.Cells.Item(6, 8).Value = 2.093154430389404
$ws.Cells.Item(6, 9).Value = 1.105769157409668

$ws.Cells.Item(7, 2).Value = 0.6395400467204037
$ws.Cells.Item(7, 3).Value = 0.5291747771945471
$ws.Cells.Item(7, 4).Value = 0.1542961280643036
$ws.Cells.Item(7, 5).Value = 0.2766794353404396
$ws.Cells.Item(7, 6).Value = 0.3989229500293732
$ws.Cells.Item(7, 7).Value = 0.2233013361692429
$ws.Cells.Item(7, 8).Value = 2.052798271179199
$ws.Cells.Item(7, 9).Value = 1.084241509437561

$ws.Cells.Item(8, 1).Value = "model_7_2_18"
$ws.Cells.Item(8, 2).Value = 0.6425806570555945
$ws.Cells.Item(8, 3).Value = 0.5377169514202708
$ws.Cells.Item(8, 4).Value = 0.1706383063085508
$ws.Cells.Item(8, 5).Value = 0.2905640041493702
$ws.Cells.Item(8, 6).Value = 0.3955579102039337
$ws.Cells.Item(8, 7).Value = 0.219249963760376
$ws.Cells.Item(8, 8).Value = 2.01313042640686
$ws.Cells.Item(8, 9).Value = 1.063428997993469

$ws.Cells.Item(9, 1).Value = "model_7_2_17"
$ws.Cells.Item(9, 2).Value = 0.6470907270654956
$ws.Cells.Item(9, 3).Value = 0.5473592131583177
$ws.Cells.Item(9, 4).Value = 0.191328695543748
$ws.Cells.Item(9, 5).Value = 0.3079456225517249
$ws.Cells.Item(9, 6).Value = 0.3905665576457977
$ws.Cells.Item(9, 7).Value = 0.2146768569946289
$ws.Cells.Item(9, 8).Value = 1.962908029556274
$ws.Cells.Item(9, 9).Value = 1.037374258041382

$ws.Cells.Item(10, 1).Value = "model_7_2_16"
$ws.Cells.Item(10, 2).Value = 0.6545934685205541
$ws.Cells.Item(10, 3).Value = 0.558818370233371
$ws.Cells.Item(10, 4).Value = 0.2219113151844232
$ws.Cells.Item(10, 5).Value = 0.3331702326863785
$ws.Cells.Item(10, 6).Value = 0.3822632133960724
$ws.Cells.Item(10, 7).Value = 0.2092420607805252
$ws.Cells.Item(10, 8).Value = 1.888674139976501
$ws.Cells.Item(10, 9).Value = 0.9995629787445068

$ws.Cells.Item(11, 1).Value = "model_7_2_15"
$ws.Cells.Item(11, 2).Value = 0.6686409342210216
$ws.Cells.Item(11, 3).Value = 0.5775400068332394
$ws.Cells.Item(11, 4).Value = 0.271831288714754
$ws.Cells.Item(11, 5).Value = 0.374346795730209
$ws.Cells.Item(11, 6).Value = 0.3667168617248535
$ws.Cells.Item(11, 7).Value = 0.2003628462553024
$ws.Cells.Item(11, 8).Value = 1.767502188682556
$ws.Cells.Item(11, 9).Value = 0.9378403425216675

$ws.Cells.Item(12, 1).Value = "model_7_2_14"
$ws.Cells.Item(12, 2).Value = 0.6792219051892892
$ws.Cells.Item(12, 3).Value = 0.6002074095102625
$ws.Cells.Item(12, 4).Value = 0.3118141047768228
$ws.Cells.Item(12, 5).Value = 0.4086118224869449
$ws.Cells.Item(12, 6).Value = 0.3550068438053131
$ws.Cells.Item(12, 7).Value = 0.1896122097969055
$ws.Cells.Item(12, 8).Value = 1.670450925827026
$ws.Cells.Item(12, 9).Value = 0.8864778280258179

$ws.Cells.Item(13, 2).Value = 0.6910822488842499
$ws.Cells.Item(13, 3).Value = 0.6139619981336619
$ws.Cells.Item(13, 4).Value = 0.3568207919521486
$ws.Cells.Item(13, 5).Value = 0.4452129471071696
$ws.Cells.Item(13, 6).Value = 0.3418809473514557
$ws.Cells.Item(13, 7).Value = 0.1830887347459793
$ws.Cells.Item(13, 8).Value = 1.56120502948761
$ws.Cells.Item(13, 9).Value = 0.8316135406494141

$ws.Cells.Item(14, 1).Value = "model_7_2_12"
$ws.Cells.Item(14, 2).Value = 0.7161230938233912
$ws.Cells.Item(14, 3).Value = 0.6315940593198217
$ws.Cells.Item(14, 4).Value = 0.4474023995817296
$ws.Cells.Item(14, 5).Value = 0.5171925263993742
$ws.Cells.Item(14, 6).Value = 0.3141681253910065
$ws.Cells.Item(14, 7).Value = 0.1747262626886368
$ws.Cells.Item(14, 8).Value = 1.341333985328674
$ws.Cells.Item(14, 9).Value = 0.7237176895141602

$ws.Cells.Item(15, 1).Value = "model_7_2_8"
$ws.Cells.Item(15, 2).Value = 0.719330732751761
$ws.Cells.Item(15, 3).Value = 0.6504892933909991
$ws.Cells.Item(15, 4).Value = 0.5390870816829414
$ws.Cells.Item(15, 5).Value = 0.5902239770280563
$ws.Cells.Item(15, 6).Value = 0.3106181919574738
$ws.Cells.Item(15, 7).Value = 0.1657647043466568
$ws.Cells.Item(15, 8).Value = 1.118785381317139
$ws.Cells.Item(15, 9).Value = 0.6142451763153076

$ws.Cells.Item(16, 1).Value = "model_7_2_7"
$ws.Cells.Item(16, 2).Value = 0.7243583226486996
$ws.Cells.Item(16, 3).Value = 0.627294750486317
$ws.Cells.Item(16, 4).Value = 0.5742183120271575
$ws.Cells.Item(16, 5).Value = 0.613109934483705
$ws.Cells.Item(16, 6).Value = 0.3050541281700134
$ws.Cells.Item(16, 7).Value = 0.1767653077840805
$ws.Cells.Item(16, 8).Value = 1.033510565757751
$ws.Cells.Item(16, 9).Value = 0.5799396634101868

$ws.Cells.Item(17, 1).Value = "model_7_2_11"
$ws.Cells.Item(17, 2).Value = 0.7259523796117096
$ws.Cells.Item(17, 3).Value = 0.6459265426592284
$ws.Cells.Item(17, 4).Value = 0.4902542567131076
$ws.Cells.Item(17, 5).Value = 0.5522475899563033
$ws.Cells.Item(17, 6).Value = 0.3032899498939514
$ws.Cells.Item(17, 7).Value = 0.1679287105798721
$ws.Cells.Item(17, 8).Value = 1.237318515777588
$ws.Cells.Item(17, 9).Value = 0.6711709499359131

$ws.Cells.Item(18, 1).Value = "model_7_2_9"
$ws.Cells.Item(18, 2).Value = 0.7272035638586224
$ws.Cells.Item(18, 3).Value = 0.658551006277297
$ws.Cells.Item(18, 4).Value = 0.5305686063390919
$ws.Cells.Item(18, 5).Value = 0.5850830452681886
$ws.Cells.Item(18, 6).Value = 0.3019053041934967
$ws.Cells.Item(18, 7).Value = 0.1619412302970886
$ws.Cells.Item(18, 8).Value = 1.13946259021759
$ws.Cells.Item(18, 9).Value = 0.6219514012336731

$ws.Cells.Item(19, 1).Value = "model_7_2_10"
$ws.Cells.Item(19, 2).Value = 0.7379779382948537
$ws.Cells.Item(19, 3).Value = 0.6688696032116639
$ws.Cells.Item(19, 4).Value = 0.5347505180107608
$ws.Cells.Item(19, 5).Value = 0.5899985442126879
$ws.Cells.Item(19, 6).Value = 0.2899812161922455
$ws.Cells.Item(19, 7).Value = 0.1570473462343216
$ws.Cells.Item(19, 8).Value = 1.129311800003052
$ws.Cells.Item(19, 9).Value = 0.6145831346511841

$ws.Cells.Item(20, 2).Value = 0.7718399217435846
$ws.Cells.Item(20, 3).Value = 0.6960843736813098
$ws.Cells.Item(20, 4).Value = 0.744843520438242
$ws.Cells.Item(20, 5).Value = 0.7546543583813384
$ws.Cells.Item(20, 6).Value = 0.2525059878826141
$ws.Cells.Item(20, 7).Value = 0.1441400349140167
$ws.Cells.Item(20, 8).Value = 0.6193476915359497
$ws.Cells.Item(20, 9).Value = 0.3677676618099213

$ws.Cells.Item(21, 2).Value = 0.776001008969595
$ws.Cells.Item(21, 3).Value = 0.7986581691314508
$ws.Cells.Item(21, 4).Value = 0.7783072570247671
$ws.Cells.Item(21, 5).Value = 0.7973368401491346
$ws.Cells.Item(21, 6).Value = 0.2479008585214615
$ws.Cells.Item(21, 7).Value = 0.0954916849732399
$ws.Cells.Item(21, 8).Value = 0.5381203293800354
$ws.Cells.Item(21, 9).Value = 0.3037875890731812

$ws.Cells.Item(22, 1).Value = "model_7_2_4"
$ws.Cells.Item(22, 2).Value = 0.7869848000232615
$ws.Cells.Item(22, 3).Value = 0.8205763323727913
$ws.Cells.Item(22, 4).Value = 0.8240879264210397
$ws.Cells.Item(22, 5).Value = 0.8358947724192284
$ws.Cells.Item(22, 6).Value = 0.2357450574636459
$ws.Cells.Item(22, 7).Value = 0.08509642630815506
$ws.Cells.Item(22, 8).Value = 0.4269957840442657
$ws.Cells.Item(22, 9).Value = 0.2459901124238968

$ws.Cells.Item(23, 1).Value = "model_7_2_0"
$ws.Cells.Item(23, 2).Value = 0.7913315794707996
$ws.Cells.Item(23, 3).Value = 0.8513022770348379
$ws.Cells.Item(23, 4).Value = 0.9145406194115412
$ws.Cells.Item(23, 5).Value = 0.9099691951498021
$ws.Cells.Item(23, 6).Value = 0.2309344559907913
$ws.Cells.Item(23, 7).Value = 0.07052382826805115
$ws.Cells.Item(23, 8).Value = 0.2074376940727234
$ws.Cells.Item(23, 9).Value = 0.1349541693925858

$ws.Cells.Item(24, 1).Value = "model_7_2_1"
$ws.Cells.Item(24, 2).Value = 0.795967551032936
$ws.Cells.Item(24, 3).Value = 0.8524914284058274
$ws.Cells.Item(24, 4).Value = 0.9045153704589464
$ws.Cells.Item(24, 5).Value = 0.9025289727722448
$ws.Cells.Item(24, 6).Value = 0.2258038073778152
$ws.Cells.Item(24, 7).Value = 0.0699598491191864
$ws.Cells.Item(24, 8).Value = 0.231772243976593
$ws.Cells.Item(24, 9).Value = 0.1461069136857986

$ws.Cells.Item(25, 1).Value = "model_7_2_3"
$ws.Cells.Item(25, 2).Value = 0.7980371216215562
$ws.Cells.Item(25, 3).Value = 0.8287629174830111
$ws.Cells.Item(25, 4).Value = 0.8766962744332921
$ws.Cells.Item(25, 5).Value = 0.8773556057548243
$ws.Cells.Item(25, 6).Value = 0.2235133945941925
$ws.Cells.Item(25, 7).Value = 0.08121372014284134
$ws.Cells.Item(25, 8).Value = 0.2992982268333435
$ws.Cells.Item(25, 9).Value = 0.1838412433862686

$ws.Cells.Item(26, 1).Value = "model_7_2_2"
$ws.Cells.Item(26, 2).Value = 0.8044162390134412
$ws.Cells.Item(26, 3).Value = 0.8491182373929974
$ws.Cells.Item(26, 4).Value = 0.9062614204759269
$ws.Cells.Item(26, 5).Value = 0.9032945784023549
$ws.Cells.Item(26, 6).Value = 0.2164535820484161
$ws.Cells.Item(26, 7).Value = 0.07155966758728027
$ws.Cells.Item(26, 8).Value = 0.2275339961051941
$ws.Cells.Item(26, 9).Value = 0.1449593007564545
